$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change the mean value in E4 from 4 to 6 (the formula in F10 will recalc automatically)
$ws.Range("E4").Value = 6

# 2. Update the caption in C9: "Value of x for which the " -> "Value of x for which the         "
#    (extra trailing spaces appended after "the")
$ws.Range("C9").Value = "Value of x for which the         "

# 3. Update the caption in C10: "probability  P{X<x} = p :   x =  " -> "probability  P{X<x} = p :     x =  "
#    (two extra spaces inserted right before "x =")
$ws.Range("C10").Value = "probability  P{X<x} = p :     x =  "

# 4. Move the active selection to G17, matching the saved selection state
[void]$ws.Range("G17").Select()
